$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Strip existing cell formatting (drops the bold header font / text-number-format
# styles) before rewriting content, so everything lands on the default style.
$ws.Range("A1:E9").ClearFormats()
$ws.Range("A2:E9").ClearContents()

# Row 2 - "non" example task
$ws.Cells.Item(2, 1).Value = 2023
$ws.Cells.Item(2, 2).Value = 23
$ws.Cells.Item(2, 3).Value = "contoh task weekly non"
$ws.Cells.Item(2, 4).Value = "NON"

# Row 3 - "result" example task
$ws.Cells.Item(3, 1).Value = 2023
$ws.Cells.Item(3, 2).Value = 23
$ws.Cells.Item(3, 3).Value = "contoh task weekly result"
$ws.Cells.Item(3, 4).Value = "RESULT"
$ws.Cells.Item(3, 5).Value = 10000

# Widen column C to fit the longer task text (~29.86 chars wide)
$ws.Columns.Item(3).ColumnWidth = 29

# Match the author's final selection / window size
$ws.Range("C3").Select()
$excel.ActiveWindow.Width = 28800
$excel.ActiveWindow.Height = 12315
